# Test_Tablas_conciliacion_PBI.xlsx - "Directorio" sheet updates
# Commit: "Se incluyen las carpetas y los cambios de la tablas_conciliacion"
#   - Row 2 (existing CO_KUSHKI_TEST entry) now lives under a "MAN" subfolder
#     and is marked Inactivo; its csv/xlsx path formulas gain a "/<Carpeta>/"
#     segment built from the (new) column C value.
#   - Two new rows (3 and 4) are added for CO_WOMPI_CONC_MAN_VEN (MAN) and
#     CO_WOMPI_CONC_SIM_VEN (SIM), following the same pattern as row 2 but
#     with H filled as a literal (not a formula referencing E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Directorio")

# --- Row 2: add the "Carpeta" (MAN) and flip Estado to Inactivo ---
$ws.Range("C2").Value = "MAN"
$ws.Range("F2").Value = "Inactivo"
$ws.Range("K2").Formula = '="C:/Users/Usuario/Desktop/git_project_1/Fuente_Power_BI/"&B2&"/"&C2&"/"&H2&".csv"'
$ws.Range("L2").Formula = '="C:/Users/Usuario/Desktop/git_project_1/Fuente_Power_BI/"&B2&"/"&C2&"/"&H2&".xlsx"'

# --- Row 3: CO_WOMPI_CONC_MAN_VEN / MAN ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "CO"
$ws.Range("C3").Value = "MAN"
$ws.Range("E3").Value = "CO_WOMPI_CONC_MAN_VEN"
$ws.Range("F3").Value = "Activo"
$ws.Range("H3").Value = "CO_WOMPI_CONC_MAN_VEN"
$ws.Range("J3").Formula = '="C:/Users/Usuario/Desktop/git_project_1/QUERY_SQL/"&B3&"/"&E3&".sql"'
$ws.Range("K3").Formula = '="C:/Users/Usuario/Desktop/git_project_1/Fuente_Power_BI/"&B3&"/"&C3&"/"&H3&".csv"'
$ws.Range("L3").Formula = '="C:/Users/Usuario/Desktop/git_project_1/Fuente_Power_BI/"&B3&"/"&C3&"/"&H3&".xlsx"'

# --- Row 4: CO_WOMPI_CONC_SIM_VEN / SIM ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "CO"
$ws.Range("C4").Value = "SIM"
$ws.Range("E4").Value = "CO_WOMPI_CONC_SIM_VEN"
$ws.Range("F4").Value = "Activo"
$ws.Range("H4").Value = "CO_WOMPI_CONC_SIM_VEN"
$ws.Range("J4").Formula = '="C:/Users/Usuario/Desktop/git_project_1/QUERY_SQL/"&B4&"/"&E4&".sql"'
$ws.Range("K4").Formula = '="C:/Users/Usuario/Desktop/git_project_1/Fuente_Power_BI/"&B4&"/"&C4&"/"&H4&".csv"'
$ws.Range("L4").Formula = '="C:/Users/Usuario/Desktop/git_project_1/Fuente_Power_BI/"&B4&"/"&C4&"/"&H4&".xlsx"'

# --- View state: selection now rests on the newly extended L column block ---
$ws.Activate()
$ws.Range("L2:L4").Select()
